$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data currently spans 2007年..2020年 (rows 2-15). The new data should
# span 2010年..2022年 (rows 2-14): drop the three oldest years (2007-2009)
# and append two new years (2021, 2022) at the end.

# Drop the 2007年 row three times in a row (each delete shifts rows up)
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()
$ws.Rows.Item(2).Delete()

# Data now occupies rows 2 (2010年) .. 12 (2020年). Append 2021年 and 2022年.
$ws.Cells.Item(13, 1).Value = "2021年"
$ws.Cells.Item(13, 2).Value = 332290550.893184
$ws.Cells.Item(13, 3).Value = 79211415.315769
$ws.Cells.Item(13, 4).Value = 449001668.215651
$ws.Cells.Item(13, 5).Value = 156086389.081515
$ws.Cells.Item(13, 6).Value = 116711117.322467
$ws.Cells.Item(13, 7).Value = 76874973.765746

$ws.Cells.Item(14, 1).Value = "2022年"
# 人身险保费收入/人身险支出 are not yet available for 2022年 -- keep as blank
# text cells (matches source workbook, which stores them as empty text cells)
# with the default, unformatted cell style.
$ws.Cells.Item(14, 2).Value = "'"
$ws.Cells.Item(14, 2).ClearFormats()
$ws.Cells.Item(14, 3).Value = "'"
$ws.Cells.Item(14, 3).ClearFormats()
$ws.Cells.Item(14, 4).Value = 469570000
$ws.Cells.Item(14, 5).Value = 154850000
$ws.Cells.Item(14, 6).Value = 127120000
$ws.Cells.Item(14, 7).Value = 77570000

# Match the bold/centered/bordered style used by the other year cells in
# column A for the two newly appended rows (copy format from A12, which
# already carries the correct look for a year label cell).
$ws.Range("A12").Copy() | Out-Null
$ws.Range("A13:A14").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
